# Actualizacion automatica 2025-11-24 08:30:09
# Updates monthly sales figures (PORCELANATO group, CERAMICAS AL COSTO / F.V - AREA ANDINA)
# across the three report sheets, keeping dependent totals/percentages in sync.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M6").Value = 415.74
$ws1.Range("M12").Value = 5890.54
$ws1.Range("M26").Value = "5 de 24"

$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F6").Value = 859.1799999999999
$ws2.Range("F12").Value = 5890.54
$ws2.Range("F26").Value = 34791.05

$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Columns.Item(5).ColumnWidth = 22.1666666666667
$ws3.Range("D12").Value = 34347.61
$ws3.Range("E12").Value = 353.3899999999994
$ws3.Range("F12").Value = 0.9898161436269848
$ws3.Range("D14").Value = 34791.05
$ws3.Range("E14").Value = 5986.69058948192
$ws3.Range("F14").Value = 0.8531872903466822
